$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson 22 (row 26): the old "Practice - recap" slot is repurposed for the
# new Java 8 Stream API lesson.
$ws.Range("C26").Value = "Java 8 (Stream API,)"

# Lesson 24 (row 27): fill in the lesson name, duration and date; row grows
# a bit taller to fit the wrapped title text.
$ws.Rows.Item(27).RowHeight = 14.9
$ws.Range("C27").Value = "Java 8 continued"
$ws.Range("D27").Value2 = 2

# Bring over the date formatting from the row above, then set the date itself
# (19-Dec-2020) without picking up a stray time-of-day component.
$ws.Range("E26").Copy($ws.Range("E27"))
$ws.Range("E27").Value2 = 44184

# Reflect where the author ended up scrolled to / selected when done editing.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C26").Select()
